$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "259.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.67%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.34%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.705"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06029"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.47%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.676"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.42%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8597"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.24%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9247"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.90%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1398"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.82%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04951"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "25.57%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.00%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03094"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09130"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.53%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006072"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006103"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.48%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.50%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.07%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.165"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.82%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.42%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.37%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.127"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.75%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04226"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.21%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.25%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004037"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-21.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03843"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.35%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1115"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.11%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004010"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01506"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "31.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002199"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.05%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005132"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.79%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05455"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-9.09%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-22.84%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
